$d = $word.ActiveDocument

# Locate the paragraph holding the "{{ messages_sentiments }}" placeholder.
foreach ($p in $d.Paragraphs) {
    $pText = $p.Range.Text
    if ($pText -like "*messages_sentiments*") {

        $pStart = $p.Range.Start
        $pTextEnd = $p.Range.End - 1   # exclude trailing paragraph mark

        # Range covering just the "messages_sentiments" token
        # (3 chars of "{{ " precede it, 3 chars of " }}" follow it).
        $token = $d.Range($pStart + 3, $pTextEnd - 3)

        # Replace the token text in place (keeps a single run for now).
        $token.Text = "image"

        # Recompute the (now shorter) token range and force it to become
        # its own run by toggling a formatting property on it and back,
        # which splits the paragraph into three runs:
        #   "{{ "  /  "image"  /  " }}"
        $newPTextEnd = $p.Range.End - 1
        $newToken = $d.Range($pStart + 3, $newPTextEnd - 3)
        $newToken.Bold = 1
        $newToken.Bold = 0

        break
    }
}
